$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values as plain text in the source data (e.g. to keep
# exact trailing zeros like "94.30" or multi-dot groupings like "41.325.36").
# Force Text format on every Price cell we touch so Excel keeps the literal string
# instead of silently re-interpreting it as a number.
$priceCells = 'D2','D3','D5','D6','D9','D10','D11','D13','D14','D15','D16','D17','D18','D20','D21','D22','D23','D27','D30','D31','D32','D35','D36','D38','D42','D44','D46','D48','D49','D50'
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '41.325.36'
$ws.Range("E2").Value = '  -3.31%  '
$ws.Range("D3").Value = '2.468.18'
$ws.Range("E3").Value = '  -2.39%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '314.39'
$ws.Range("E5").Value = '  +1.65%  '
$ws.Range("D6").Value = '94.30'
$ws.Range("E6").Value = '  -7.39%  '
$ws.Range("E7").Value = '  -2.80%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.502'
$ws.Range("E9").Value = '  -4.78%  '
$ws.Range("D10").Value = '33.56'
$ws.Range("E10").Value = '  -6.68%  '
$ws.Range("D11").Value = '0.0780'
$ws.Range("E11").Value = '  -3.11%  '
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").Value = '7.01'
$ws.Range("E13").Value = '  -4.38%  '
$ws.Range("D14").Value = '2.850.67'
$ws.Range("E14").Value = '  -2.18%  '
$ws.Range("D15").Value = '2.480.78'
$ws.Range("E15").Value = '  -4.13%  '
$ws.Range("D16").Value = '14.66'
$ws.Range("E16").Value = '  -6.75%  '
$ws.Range("D17").Value = '0.785'
$ws.Range("E17").Value = '  -3.01%  '
$ws.Range("D18").Value = '41.303.28'
$ws.Range("E18").Value = '  -3.28%  '
$ws.Range("E19").Value = '  -6.70%  '
$ws.Range("D20").Value = '0.0₃0920'
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("D21").Value = '11.48'
$ws.Range("E21").Value = '  -5.07%  '
$ws.Range("D22").Value = '68.01'
$ws.Range("E22").Value = '  -1.96%  '
$ws.Range("D23").Value = '237.85'
$ws.Range("E23").Value = '  -2.65%  '
$ws.Range("E24").Value = '  -3.13%  '
$ws.Range("E25").Value = '  -5.32%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = '24.45'
$ws.Range("E27").Value = '  -6.60%  '
$ws.Range("E29").Value = '  -4.43%  '
$ws.Range("D30").Value = '35.99'
$ws.Range("E30").Value = '  -8.05%  '
$ws.Range("D31").Value = '152.55'
$ws.Range("E31").Value = '  -3.11%  '
$ws.Range("D32").Value = '5.61'
$ws.Range("E32").Value = '  -3.24%  '
$ws.Range("E33").Value = '  -6.47%  '
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("D35").Value = '0.0752'
$ws.Range("E35").Value = '  -4.78%  '
$ws.Range("D36").Value = '3.02'
$ws.Range("E36").Value = '  -5.40%  '
$ws.Range("E37").Value = '  -7.25%  '
$ws.Range("D38").Value = '16.89'
$ws.Range("E38").Value = '  -6.86%  '
$ws.Range("E39").Value = '  -7.11%  '
$ws.Range("E40").Value = '  -3.83%  '
$ws.Range("E41").Value = '  +2.49%  '
$ws.Range("D42").Value = '21.05'
$ws.Range("E42").Value = '  -4.62%  '
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").Value = '1.988.92'
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("E45").Value = '  -5.07%  '
$ws.Range("D46").Value = '3.07'
$ws.Range("E46").Value = '  -6.73%  '
$ws.Range("E47").Value = '  -1.78%  '
$ws.Range("D48").Value = '69.84'
$ws.Range("E48").Value = '  -3.23%  '
$ws.Range("D49").Value = '75.83'
$ws.Range("E49").Value = '  -5.64%  '
$ws.Range("D50").Value = '96.89'
$ws.Range("E50").Value = '  -4.38%  '
$ws.Range("E51").Value = '  -6.27%  '
